# "Forms synchronized with database" -- mandatory_fields.xlsx / Sheet1
# Re-labels the three header columns and reshuffles which mandatory-field
# name appears in which row/column so the sheet matches the DB schema.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) ---------------------------------------------------
$ws.Range("A1").Value = "Need support and emergency"
$ws.Range("B1").Value = "Offer support"
$ws.Range("C1").Value = "Message"

# ---- Rows 2-3 unchanged (message_id / message_type) -----------------------

# ---- Row 4: title stays in all three columns, but the border style used
# on column C moves to column B (swap the two cell formats). -----------------
$ws.Range("C4").Copy()
$ws.Range("Z1").PasteSpecial(-4122)
$ws.Range("B4").Copy()
$ws.Range("C4").PasteSpecial(-4122)
$ws.Range("Z1").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("Z1").Clear()

# ---- Row 5-6: shift field names up one slot --------------------------------
$ws.Range("A5").Value = "time_start"
$ws.Range("A6").Value = "time_stop"

# ---- Row 7: date_of_creation / date of creation swap across columns -------
$ws.Range("A7").Value = "date_of_creation"
$ws.Range("B7").Value = "date of creation"
$ws.Range("C7").Value = "date_of_creation"

# ---- Row 8: description takes over A and C ---------------------------------
$ws.Range("A8").Value = "description"
$ws.Range("C8").Value = "description"

# ---- Row 9: category moves to A and (new) B; C9 becomes blank -------------
$ws.Range("A9").Value = "category"
$ws.Range("A9").Copy()
$ws.Range("B9").PasteSpecial(-4122)
$ws.Range("B9").Value = "category"
$ws.Range("C9").ClearContents()

# ---- Row 10: person_contact moves to A and B; C10 becomes blank -----------
$ws.Range("A10").Value = "person_contact"
$ws.Range("B10").Value = "person_contact"
$ws.Range("C10").ClearContents()

# ---- Row 11: location moves to A; B/C stay blank ---------------------------
$ws.Range("A11").Value = "location"

# ---- Row 12: A12 is dropped entirely (B12/C12 remain as blank cells) ------
$ws.Range("A12").Clear()

# ---- Column widths (A wider, B narrower, C/D set to match) ----------------
$ws.Columns.Item(1).ColumnWidth = 28.833333333333336
$ws.Columns.Item(2).ColumnWidth = 21.666666666666668
$ws.Columns.Item(3).ColumnWidth = 26.666666666666668
$ws.Columns.Item(4).ColumnWidth = 26.666666666666668

# ---- Selection moves from C9 to C7 -----------------------------------------
$ws.Range("C7").Select() | Out-Null
